# implement batch import of controls
#
# The reagents import template gains two new columns so a batch import can
# also carry control metadata: "manufacturer" (D) and "supplier" (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1). These become new shared-string entries.
$ws.Range("D1").Value = "manufacturer"
$ws.Range("E1").Value = "supplier"

# Size the new "manufacturer" column to fit its header text, same as Excel's
# own best-fit column sizing would do after typing a header in.
$ws.Columns.Item(4).ColumnWidth = 11.5

# Leave the selection on the newly added cell, matching the saved view state.
$ws.Range("D2").Select() | Out-Null
